# Insert a new weekly record as row 153 in the "Hortaliza, Feria Lagunitas de
# Puerto Montt - Zapallo" sheet. Existing rows 153:178 shift down to 154:179
# (dimension grows from A1:R178 to A1:R179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(153).Insert()

$ws.Range("A153").Value = 4
$ws.Range("B153").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C153").Value = "Los Lagos"
$ws.Range("D153").Value = 44476
$ws.Range("E153").Value = 10
$ws.Range("F153").Value = 100112045
$ws.Range("G153").Value = "Zapallo"
$ws.Range("H153").Value = "Paine"
$ws.Range("I153").Value = "1a (guarda)"
$ws.Range("J153").Value = 500
$ws.Range("K153").Value = 600
$ws.Range("L153").Value = 600
$ws.Range("M153").Value = 600
$ws.Range("N153").Value = "$/kilo (volumen en unidades)"
$ws.Range("O153").Value = "Región Metropolitana"
$ws.Range("P153").Value = 600
$ws.Range("Q153").Value = 1
$ws.Range("R153").Value = "Hortaliza"
